$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinates in row 3 to whole numbers
$ws.Range("Q3").Value = 772246
$ws.Range("R3").Value = 7120217

# Remove the start/end time cells for row 3 (they become empty cells entirely)
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
